$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I8").Value = 'aa'
$ws.Range("J8").Value = 'Agree/Accept'
$ws.Range("I17").Value = 'aa'
$ws.Range("J17").Value = 'Agree/Accept'
$ws.Range("I22").Value = '%'
$ws.Range("J22").Value = 'Uninterpretable'
$ws.Range("I24").Value = 'b'
$ws.Range("J24").Value = 'Acknowledge (Backchannel)'
$ws.Range("I30").Value = 'sd'
$ws.Range("J30").Value = 'Statement-non-opinion'
$ws.Range("I35").Value = 'aa'
$ws.Range("J35").Value = 'Agree/Accept'
$ws.Range("I50").Value = 'aa'
$ws.Range("J50").Value = 'Agree/Accept'
$ws.Range("I56").Value = 'sv'
$ws.Range("J56").Value = 'Statement-opinion'
$ws.Range("I63").Value = 'sd'
$ws.Range("J63").Value = 'Statement-non-opinion'
$ws.Range("I72").Value = 'aa'
$ws.Range("J72").Value = 'Agree/Accept'
$ws.Range("I74").Value = 'ba'
$ws.Range("J74").Value = 'Appreciation'
$ws.Range("I78").Value = 'ba'
$ws.Range("J78").Value = 'Appreciation'
$ws.Range("I101").Value = 'b'
$ws.Range("J101").Value = 'Acknowledge (Backchannel)'
$ws.Range("I105").Value = 'ba'
$ws.Range("J105").Value = 'Appreciation'
$ws.Range("I106").Value = 'b'
$ws.Range("J106").Value = 'Acknowledge (Backchannel)'
$ws.Range("I110").Value = 'ba'
$ws.Range("J110").Value = 'Appreciation'
$ws.Range("I113").Value = 'ba'
$ws.Range("J113").Value = 'Appreciation'
$ws.Range("I115").Value = 'ba'
$ws.Range("J115").Value = 'Appreciation'
$ws.Range("I129").Value = 'ba'
$ws.Range("J129").Value = 'Appreciation'
$ws.Range("I140").Value = 'sd'
$ws.Range("J140").Value = 'Statement-non-opinion'
$ws.Range("I151").Value = 'sd'
$ws.Range("J151").Value = 'Statement-non-opinion'
$ws.Range("I191").Value = 'b'
$ws.Range("J191").Value = 'Acknowledge (Backchannel)'
$ws.Range("I201").Value = 'ba'
$ws.Range("J201").Value = 'Appreciation'
$ws.Range("I220").Value = 'ba'
$ws.Range("J220").Value = 'Appreciation'
$ws.Range("I222").Value = 'ba'
$ws.Range("J222").Value = 'Appreciation'
$ws.Range("I227").Value = 'sd'
$ws.Range("J227").Value = 'Statement-non-opinion'
$ws.Range("I230").Value = 'b'
$ws.Range("J230").Value = 'Acknowledge (Backchannel)'
$ws.Range("I256").Value = 'sd'
$ws.Range("J256").Value = 'Statement-non-opinion'
$ws.Range("I257").Value = 'sd'
$ws.Range("J257").Value = 'Statement-non-opinion'
$ws.Range("I259").Value = 'sd'
$ws.Range("J259").Value = 'Statement-non-opinion'
$ws.Range("I266").Value = 'ba'
$ws.Range("J266").Value = 'Appreciation'
$ws.Range("I270").Value = 'b'
$ws.Range("J270").Value = 'Acknowledge (Backchannel)'
$ws.Range("I273").Value = 'sd'
$ws.Range("J273").Value = 'Statement-non-opinion'
$ws.Range("I281").Value = 'sd'
$ws.Range("J281").Value = 'Statement-non-opinion'
$ws.Range("I284").Value = 'b'
$ws.Range("J284").Value = 'Acknowledge (Backchannel)'
$ws.Range("I288").Value = 'sd'
$ws.Range("J288").Value = 'Statement-non-opinion'
$ws.Range("I292").Value = 'sd'
$ws.Range("J292").Value = 'Statement-non-opinion'
$ws.Range("I293").Value = 'sd'
$ws.Range("J293").Value = 'Statement-non-opinion'
$ws.Range("I295").Value = 'b'
$ws.Range("J295").Value = 'Acknowledge (Backchannel)'
$ws.Range("I296").Value = 'aa'
$ws.Range("J296").Value = 'Agree/Accept'
$ws.Range("I301").Value = 'sd'
$ws.Range("J301").Value = 'Statement-non-opinion'
$ws.Range("I303").Value = 'sv'
$ws.Range("J303").Value = 'Statement-opinion'
$ws.Range("I306").Value = 'b'
$ws.Range("J306").Value = 'Acknowledge (Backchannel)'
$ws.Range("I310").Value = 'aa'
$ws.Range("J310").Value = 'Agree/Accept'
$ws.Range("I311").Value = 'b'
$ws.Range("J311").Value = 'Acknowledge (Backchannel)'
$ws.Range("I316").Value = 'sv'
$ws.Range("J316").Value = 'Statement-opinion'
$ws.Range("I322").Value = 'sv'
$ws.Range("J322").Value = 'Statement-opinion'
$ws.Range("I323").Value = 'sd'
$ws.Range("J323").Value = 'Statement-non-opinion'
$ws.Range("I326").Value = 'b'
$ws.Range("J326").Value = 'Acknowledge (Backchannel)'
$ws.Range("I332").Value = 'sd'
$ws.Range("J332").Value = 'Statement-non-opinion'
$ws.Range("I343").Value = '%'
$ws.Range("J343").Value = 'Uninterpretable'
$ws.Range("I351").Value = 'sd'
$ws.Range("J351").Value = 'Statement-non-opinion'
$ws.Range("I358").Value = 'sd'
$ws.Range("J358").Value = 'Statement-non-opinion'
$ws.Range("I359").Value = 'sd'
$ws.Range("J359").Value = 'Statement-non-opinion'
$ws.Range("I361").Value = 'sv'
$ws.Range("J361").Value = 'Statement-opinion'
$ws.Range("I364").Value = 'b'
$ws.Range("J364").Value = 'Acknowledge (Backchannel)'
$ws.Range("I380").Value = 'b'
$ws.Range("J380").Value = 'Acknowledge (Backchannel)'
$ws.Range("I382").Value = 'ba'
$ws.Range("J382").Value = 'Appreciation'
$ws.Range("I386").Value = 'sd'
$ws.Range("J386").Value = 'Statement-non-opinion'
$ws.Range("I387").Value = 'sd'
$ws.Range("J387").Value = 'Statement-non-opinion'
$ws.Range("I429").Value = 'b'
$ws.Range("J429").Value = 'Acknowledge (Backchannel)'
$ws.Range("I440").Value = 'sd'
$ws.Range("J440").Value = 'Statement-non-opinion'
$ws.Range("I458").Value = '%'
$ws.Range("J458").Value = 'Uninterpretable'
$ws.Range("I463").Value = 'b'
$ws.Range("J463").Value = 'Acknowledge (Backchannel)'
$ws.Range("I465").Value = 'sd'
$ws.Range("J465").Value = 'Statement-non-opinion'
$ws.Range("I497").Value = 'ba'
$ws.Range("J497").Value = 'Appreciation'
$ws.Range("I501").Value = 'aa'
$ws.Range("J501").Value = 'Agree/Accept'
$ws.Range("I508").Value = 'sd'
$ws.Range("J508").Value = 'Statement-non-opinion'
$ws.Range("I512").Value = 'ba'
$ws.Range("J512").Value = 'Appreciation'
$ws.Range("I521").Value = 'sd'
$ws.Range("J521").Value = 'Statement-non-opinion'
$ws.Range("I538").Value = 'ba'
$ws.Range("J538").Value = 'Appreciation'
$ws.Range("I540").Value = 'b'
$ws.Range("J540").Value = 'Acknowledge (Backchannel)'
